# Daily attendance processing - 2026-01-15 08:45:27
#
# Applies the scraped changes to "Session Analysis Results":
#   1. Summary counters L7 (Missing Sessions) and L8 (Pending Sessions).
#   2. "Recorded By" text in column G: reorder "System, <email>" -> "<email>, System".
#   3. Swap the "Missing"/"Pending" pair in columns P/Q for the B1A1 weekly summary rows.
#   4. Six newly-"Not Recorded" placeholder rows (B1D1/B1D2/B1E1/B1E2/B1F1/B1F2 for
#      15/01/2026): recolor from the "Pending" (yellow) style to the "Not Recorded"
#      (pink) style and update the Status text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Summary counters (K7/L7 = Missing Sessions, K8/L8 = Pending Sessions)
$ws.Range("L7").Value2 = 33
$ws.Range("L8").Value2 = 36

# 2. "Recorded By" column: swap the order of the two names/addresses.
$recordedByRows = @(8, 9, 10, 12, 14, 15, 17, 18, 34, 35, 36, 38, 40, 41, 43, 44, 60, 61, 62, 64, 66, 67, 69, 70, 86, 87, 88, 90, 92, 93, 95, 96, 112, 113, 114, 116, 118, 119, 121, 122, 138, 139, 140, 142, 144, 145, 147, 148, 164, 167, 170, 174, 191, 194, 197, 201, 218, 221, 224, 228, 245, 248, 251, 255, 272, 275, 278, 282, 299, 302, 305, 309)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value2 = "dnasr281@gmail.com, System"
}

# 3. Swap columns P/Q (Missing / Pending weekly counts) for rows 21-23, 25-26;
#    row 24 shifts one Pending into Missing (4/4 -> 5/3) rather than a pure swap.
$swapRows = @(21, 22, 23, 25, 26)
foreach ($r in $swapRows) {
    $pCell = $ws.Range("P$r")
    $qCell = $ws.Range("Q$r")
    $pVal = $pCell.Value2
    $qVal = $qCell.Value2
    $pCell.Value2 = $qVal
    $qCell.Value2 = $pVal
}
$ws.Range("P24").Value2 = 5
$ws.Range("Q24").Value2 = 3

# 4. Recolor the six placeholder rows from "Pending" style to "Not Recorded" style
#    by copying the formatting of an existing "Not Recorded" row (row 3), then fix
#    up the Status text.
$notRecordedRows = @(181, 208, 235, 262, 289, 316)
$ws.Range("A3:I3").Copy()
foreach ($r in $notRecordedRows) {
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value2 = "Not Recorded"
}
$excel.CutCopyMode = $false
